$wb = $excel.ActiveWorkbook

# New row (81) data for each of the four worksheets, in sheet order.
$rowsData = @(
    @{ A = 45867.46590277777; B = "0x01,0x90"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"; D = "0x01,0x34"; E = "0x07"; F = 400; G = [double]"5.68631262647113e+23"; H = 308; I = 7  },
    @{ A = 45867.46590277777; B = "0x01,0x7c"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"; D = "0x01,0x30"; E = "0x19"; F = 380; G = [double]"5.68432987514711e+23"; H = 304; I = 25 },
    @{ A = 45867.46590277777; B = "0x00,0x6e"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"; D = "0x00,0x60"; E = "0x15"; F = 110; G = [double]"5.68631262647113e+23"; H = 96;  I = 15 },
    @{ A = 45867.46590277777; B = "0x00,0x82"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"; D = "0x00,0x76"; E = "0x9";  F = 130; G = [double]"5.68631262647113e+23"; H = 118; I = 9  }
)

$newRow = 81

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $data = $rowsData[$i - 1]

    # Column A: numeric date/time serial, formatted like the rows above it.
    $cellA = $ws.Cells.Item($newRow, 1)
    $cellA.Value = $data.A
    $cellA.NumberFormat = $ws.Cells.Item($newRow - 1, 1).NumberFormat

    # Columns B-E: inline strings holding comma separated hex byte lists.
    $ws.Cells.Item($newRow, 2).Value = $data.B
    $ws.Cells.Item($newRow, 3).Value = $data.C
    $ws.Cells.Item($newRow, 4).Value = $data.D
    $ws.Cells.Item($newRow, 5).Value = $data.E

    # Columns F-I: plain numbers.
    $ws.Cells.Item($newRow, 6).Value = $data.F
    $ws.Cells.Item($newRow, 7).Value = $data.G
    $ws.Cells.Item($newRow, 8).Value = $data.H
    $ws.Cells.Item($newRow, 9).Value = $data.I
}
